# regen sval data to filter save games
# Replaces the computed statistic columns (B:G) for rows 2-15 with the
# regenerated values (column A dates / F win-flags stay the same; column
# G is the row sum of B:E).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
  @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 1, 6.189590430959694),
  @(0.6606524410359556, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 0, 2.960089034096801),
  @(0.6606524410359556, 1.655778082260271, 0.7527432677738641, 10.19245300693656, 1, 13.26162679800665),
  @(3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 0, 5.586269137925634),
  @(0.1190320826869504, 0.04071648406533734, 0.7527432677738641, 0.4942365360607697, 1, 1.406728370586922),
  @(1.455362044514542, 0.306821227259698, 0.1494219747398047, 0.4942365360607697, 0, 2.405841782574814),
  @(1.455362044514542, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 1, 4.358119930609447),
  @(0.6606524410359556, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 0, 3.56341032713086),
  @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 1, 6.189590430959694),
  @(3.286832544864788, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 0, 8.974608811992548),
  @(1.455362044514542, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 1, 4.358119930609447),
  @(0.1190320826869504, 0.306821227259698, 0.7527432677738641, 0.4942365360607697, 1, 1.672833113781282),
  @(3.286832544864788, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 1, 8.974608811992548),
  @(0.6606524410359556, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 0, 6.348428708163715)
)

$startRow = 2
for ($i = 0; $i -lt $newValues.Count; $i++) {
    $row = $startRow + $i
    $vals = $newValues[$i]
    for ($j = 0; $j -lt $vals.Count; $j++) {
        $col = 2 + $j   # column B = 2
        $ws.Cells.Item($row, $col).Value = $vals[$j]
    }
}
